$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "Noom-157"
$ws.Range("B5").Value = "'2023-04-26"
$ws.Range("C5").Value = "01:04:13"

# Row 6
$ws.Range("A6").Value = "Noom-157"
$ws.Range("B6").Value = "'2023-04-26"
$ws.Range("C6").Value = "01:10:14"
